# Applies "Database improvements + RES correction + Improved Outages"
# to the Annual_Generation_Statistics workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generation")

# --- RES correction: update Biofuels (column E) figures for a few countries ---
$ws.Range("E2").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("E10").Value = 6
$ws.Range("E20").Value = 25

# --- Reflect the last active cell/selection used while reviewing the sheet ---
$ws.Range("E21").Select()
